# Auto-generated edit script: updates Leve profit-tracking numbers
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets (scheduled price refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 69.61
$ws.Range("I15").Value = 69.61
$ws.Range("K15").Value = 208.83
$ws.Range("M15").Value = -39.82999999999998
$ws.Range("H38").Value = 216.625
$ws.Range("I38").Value = 104.71429
$ws.Range("J38").Value = 1000
$ws.Range("K38").Value = 314.14287
$ws.Range("L38").Value = 3000
$ws.Range("M38").Value = 57.85712999999998
$ws.Range("N38").Value = -3744
$ws.Range("H113").Value = 3066.7837
$ws.Range("I113").Value = 2968.8333
$ws.Range("J113").Value = 3486.5715
$ws.Range("K113").Value = 2968.8333
$ws.Range("L113").Value = 3486.5715
$ws.Range("M113").Value = 285.1667000000002
$ws.Range("N113").Value = -9994.5715
$ws.Range("H116").Value = 4100.654
$ws.Range("I116").Value = 4259.3687
$ws.Range("K116").Value = 4259.3687
$ws.Range("M116").Value = -817.3687
$ws.Range("H129").Value = 975763.75
$ws.Range("J129").Value = 1090502.5
$ws.Range("L129").Value = 3271507.5
$ws.Range("N129").Value = -3281507.5
$ws.Range("H132").Value = 2566353.8
$ws.Range("I132").Value = 2096.4138
$ws.Range("J132").Value = 10002700
$ws.Range("K132").Value = 6289.241399999999
$ws.Range("L132").Value = 30008100
$ws.Range("M132").Value = -3759.241399999999
$ws.Range("N132").Value = -30013160
$ws.Range("H135").Value = 1498.2
$ws.Range("I135").Value = 1498.2
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 13483.8
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -10948.8
$ws.Range("N135").Value = $null
$ws.Range("H137").Value = 985.0345
$ws.Range("I137").Value = 898.4761999999999
$ws.Range("J137").Value = 1212.25
$ws.Range("K137").Value = 2695.4286
$ws.Range("L137").Value = 3636.75
$ws.Range("M137").Value = -145.4285999999997
$ws.Range("N137").Value = -8736.75
$ws.Range("H138").Value = 3289.74
$ws.Range("I138").Value = 1440.5927
$ws.Range("J138").Value = 3973.6711
$ws.Range("K138").Value = 4321.7781
$ws.Range("L138").Value = 11921.0133
$ws.Range("M138").Value = 818.2219000000005
$ws.Range("N138").Value = -22201.0133
$ws.Range("H141").Value = 934
$ws.Range("I141").Value = 929.6
$ws.Range("K141").Value = 2788.8
$ws.Range("M141").Value = 2391.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4026.49
$ws.Range("I32").Value = 4026.49
$ws.Range("K32").Value = 4026.49
$ws.Range("M32").Value = -3739.49
$ws.Range("H61").Value = 3758.923
$ws.Range("I61").Value = 3805.3333
$ws.Range("J61").Value = 3654.5
$ws.Range("K61").Value = 3805.3333
$ws.Range("L61").Value = 3654.5
$ws.Range("M61").Value = -3593.3333
$ws.Range("N61").Value = -4078.5
$ws.Range("H74").Value = 1306.9395
$ws.Range("I74").Value = 1342.875
$ws.Range("J74").Value = 1211.1111
$ws.Range("K74").Value = 1342.875
$ws.Range("L74").Value = 1211.1111
$ws.Range("M74").Value = -468.875
$ws.Range("N74").Value = -2959.1111
$ws.Range("H77").Value = 1306.9395
$ws.Range("I77").Value = 1342.875
$ws.Range("J77").Value = 1211.1111
$ws.Range("K77").Value = 6714.375
$ws.Range("L77").Value = 6055.5555
$ws.Range("M77").Value = -2346.375
$ws.Range("N77").Value = -14791.5555
$ws.Range("H136").Value = 3758.923
$ws.Range("I136").Value = 3805.3333
$ws.Range("J136").Value = 3654.5
$ws.Range("K136").Value = 11415.9999
$ws.Range("L136").Value = 10963.5
$ws.Range("M136").Value = -8865.999899999999
$ws.Range("N136").Value = -16063.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 140
$ws.Range("I7").Value = 100
$ws.Range("K7").Value = 100
$ws.Range("M7").Value = 13
$ws.Range("H31").Value = 47614.74
$ws.Range("I31").Value = 4654.1816
$ws.Range("J31").Value = 86995.25
$ws.Range("K31").Value = 4654.1816
$ws.Range("L31").Value = 86995.25
$ws.Range("M31").Value = -4359.1816
$ws.Range("N31").Value = -87585.25
$ws.Range("H34").Value = 47614.74
$ws.Range("I34").Value = 4654.1816
$ws.Range("J34").Value = 86995.25
$ws.Range("K34").Value = 4654.1816
$ws.Range("L34").Value = 86995.25
$ws.Range("M34").Value = -4452.1816
$ws.Range("N34").Value = -87399.25
$ws.Range("H132").Value = 1324.8572
$ws.Range("I132").Value = 1046.8695
$ws.Range("J132").Value = 2603.6
$ws.Range("K132").Value = 3140.6085
$ws.Range("L132").Value = 7810.799999999999
$ws.Range("M132").Value = -610.6085000000003
$ws.Range("N132").Value = -12870.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 160
$ws.Range("I10").Value = 160
$ws.Range("K10").Value = 480
$ws.Range("M10").Value = -341
$ws.Range("H11").Value = 57338.57
$ws.Range("J11").Value = 900
$ws.Range("L11").Value = 2700
$ws.Range("N11").Value = -2980
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").Value = $null
$ws.Range("H29").Value = 100169.2
$ws.Range("J29").Value = 125186.5
$ws.Range("L29").Value = 375559.5
$ws.Range("N29").Value = -376113.5
$ws.Range("H50").Value = 101.42857
$ws.Range("I50").Value = 62
$ws.Range("J50").Value = 200
$ws.Range("K50").Value = 186
$ws.Range("L50").Value = 600
$ws.Range("M50").Value = 295
$ws.Range("N50").Value = -1562
$ws.Range("H53").Value = 101.42857
$ws.Range("I53").Value = 62
$ws.Range("J53").Value = 200
$ws.Range("K53").Value = 186
$ws.Range("L53").Value = 600
$ws.Range("M53").Value = 295
$ws.Range("N53").Value = -1562
$ws.Range("H69").Value = 462.6
$ws.Range("J69").Value = 700
$ws.Range("L69").Value = 2100
$ws.Range("N69").Value = -3722
$ws.Range("H72").Value = 462.6
$ws.Range("J72").Value = 700
$ws.Range("L72").Value = 6300
$ws.Range("N72").Value = -14412
$ws.Range("H109").Value = 25001692
$ws.Range("I109").Value = 55555936
$ws.Range("J109").Value = 2763.6365
$ws.Range("K109").Value = 166667808
$ws.Range("L109").Value = 8290.9095
$ws.Range("M109").Value = -166666768
$ws.Range("N109").Value = -10370.9095
$ws.Range("H119").Value = 6266.3335
$ws.Range("I119").Value = 1149.5
$ws.Range("J119").Value = 16500
$ws.Range("K119").Value = 3448.5
$ws.Range("L119").Value = 49500
$ws.Range("M119").Value = 1389.5
$ws.Range("N119").Value = -59176
$ws.Range("H122").Value = 23410.756
$ws.Range("I122").Value = 571.7778
$ws.Range("J122").Value = 29120.5
$ws.Range("K122").Value = 5146.000199999999
$ws.Range("L122").Value = 262084.5
$ws.Range("M122").Value = -2696.000199999999
$ws.Range("N122").Value = -266984.5
$ws.Range("H131").Value = 15213270
$ws.Range("J131").Value = 47007.41
$ws.Range("L131").Value = 141022.23
$ws.Range("N131").Value = -151102.23

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 100004
$ws.Range("J12").Value = 100004
$ws.Range("L12").Value = 100004
$ws.Range("N12").Value = -100284
$ws.Range("H122").Value = 1462849.8
$ws.Range("I122").Value = 1645580
$ws.Range("J122").Value = 1008
$ws.Range("K122").Value = 4936740
$ws.Range("L122").Value = 3024
$ws.Range("M122").Value = -4934290
$ws.Range("N122").Value = -7924
$ws.Range("H132").Value = 2239.9058
$ws.Range("I132").Value = 1818.2559
$ws.Range("J132").Value = 4053
$ws.Range("K132").Value = 5454.7677
$ws.Range("L132").Value = 12159
$ws.Range("M132").Value = -2924.7677
$ws.Range("N132").Value = -17219

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H119").Value = 48966.668
$ws.Range("J119").Value = 48966.668
$ws.Range("L119").Value = 48966.668
$ws.Range("N119").Value = -58642.668
$ws.Range("H132").Value = 2928.2
$ws.Range("I132").Value = 2021.5385
$ws.Range("J132").Value = 5138.1875
$ws.Range("K132").Value = 6064.6155
$ws.Range("L132").Value = 15414.5625
$ws.Range("M132").Value = -3534.6155
$ws.Range("N132").Value = -20474.5625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").Value = $null
$ws.Range("H136").Value = 895.55884
$ws.Range("I136").Value = 672.9583
$ws.Range("J136").Value = 1429.8
$ws.Range("K136").Value = 2018.8749
$ws.Range("L136").Value = 4289.4
$ws.Range("M136").Value = 531.1251
$ws.Range("N136").Value = -9389.4

Write-Output "applied changes"
